# "added meta data to models so easier to understand"
#
# Adds a "Meta data" section below the existing accuracy/precision/recall
# tables on Sheet1 (rows 9-15), pairing each variable code used above with
# a human-readable description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header, bold like the "Variables in model" header in B2.
$ws.Range("B9").Value = "Meta data"
$ws.Range("B9").Font.Bold = $true

# Variable codes (column C).
$ws.Range("C10").Value = "P8MEA"
$ws.Range("C12").Value = "KS2APS"
$ws.Range("C11").Value = "P8MEA_17"
$ws.Range("C13").Value = "ATT8SCR_17"
$ws.Range("C14").Value = "PTEBACHUM_E_PTQ_EE"
$ws.Range("C15").Value = "TEBACLAN_E_PTQ_EE"

# Matching descriptions (column B).
$ws.Range("B10").Value = "Progress 8 measure after adjustment for extreme scores (2019)"
$ws.Range("B12").Value = "Key stage 2 Average Points Score of the cohort at the end of key stage 4"
$ws.Range("B11").Value = "Progress 8 measure after adjustment for extreme scores (2017)"
$ws.Range("B13").Value = "Attainment 8 score (2017)"
$ws.Range("B14").Value = "% of pupils entering the English Baccalaureate Humanities subject area"
$ws.Range("B15").Value = "% of pupils entering the English Baccalaureate Language subject area"

# Leave the cursor on the last entry, matching where the author ended up.
$ws.Range("B15").Select()
